$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "rxxx"
$ws.Range("B4").Value = "cameron"
$ws.Range("C4").Value = "yet another test"
$ws.Range("D4").Value = "2025-09-27 00:54:35"
